$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 50 (shifts existing rows 50+ down by one),
# inheriting formatting from the row above as Excel normally does.
$ws.Rows("50:50").Insert()

# Fill in the new LeetCode entry: 219. Contains Duplicate II
$ws.Range("A50").Value = 219
$ws.Range("B50").Value = "Contains Duplicate II"
$ws.Range("C50").Value = "Easy"
$ws.Range("D50").Value = "Arrays,sliding window,hashmap"
$ws.Range("E50").Value = 45872
$ws.Range("F50").Value = "Python"

# The hyperlink on the "127. Word Ladder" row does not automatically
# follow the row shift, so move it from its old anchor to the new one.
$ws.Range("B128").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B129"), "https://leetcode.com/problems/word-ladder/", [System.Reflection.Missing]::Value, "https://leetcode.com/problems/word-ladder/", "127. Word Ladder")

# Update the view so the selection matches the edited area.
$ws.Range("F51").Select()
